$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(5)
$tr = $sh.TextFrame.TextRange
Write-Host "Length:" $tr.Length
Write-Host "Text:" $tr.Text
$chars = $tr.Characters(1, $tr.Length)
Write-Host "Chars text:" $chars.Text
$chars.Text = "11/9/2020"
Write-Host "After Chars text:" $tr.Text
